$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "贵州茅台"
$ws.Range("C2").Value = "利欧股份"
$ws.Range("A3").Value = "白银有色"
$ws.Range("B3").Value = "贵州茅台"
$ws.Range("C3").Value = "贵州茅台"
$ws.Range("C4").Value = "五粮液"
$ws.Range("A5").Value = "利欧股份"
$ws.Range("C5").Value = "白银有色"
$ws.Range("A6").Value = "五 粮 液"
$ws.Range("B6").Value = "蓝色光标"
$ws.Range("C6").Value = "铜陵有色"
$ws.Range("A7").Value = "浙文互联"
$ws.Range("B7").Value = "湖南白银"
$ws.Range("C7").Value = "蓝色光标"
$ws.Range("A8").Value = "中国铝业"
$ws.Range("B8").Value = "中国铝业"
$ws.Range("C8").Value = "湖南白银"
$ws.Range("A9").Value = "中国黄金"
$ws.Range("B9").Value = "紫金矿业"
$ws.Range("C9").Value = "洲际油气"
$ws.Range("A10").Value = "铜陵有色"
$ws.Range("B10").Value = "中国黄金"
$ws.Range("C10").Value = "世纪华通"
$ws.Range("A11").Value = "湖南白银"
$ws.Range("B11").Value = "五 粮 液"
$ws.Range("C11").Value = "浙文互联"
$ws.Range("A12").Value = "洲际油气"
$ws.Range("B12").Value = "浙文互联"
$ws.Range("C12").Value = "红宝丽"
$ws.Range("A13").Value = "紫金矿业"
$ws.Range("B13").Value = "湖南黄金"
$ws.Range("C13").Value = "中国铝业"
$ws.Range("A14").Value = "湖南黄金"
$ws.Range("B14").Value = "洲际油气"
$ws.Range("A15").Value = "洛阳钼业"
$ws.Range("B15").Value = "北方铜业"
$ws.Range("C15").Value = "紫金矿业"
$ws.Range("A16").Value = "晓程科技"
$ws.Range("B16").Value = "云南铜业"
$ws.Range("C16").Value = "晓程科技"
$ws.Range("A17").Value = "网宿科技"
$ws.Range("B17").Value = "洛阳钼业"
$ws.Range("C17").Value = "中国黄金"
$ws.Range("A18").Value = "江西铜业"
$ws.Range("B18").Value = "北方稀土"
$ws.Range("C18").Value = "通源石油"
$ws.Range("A19").Value = "通源石油"
$ws.Range("B19").Value = "江西铜业"
$ws.Range("C19").Value = "江西铜业"
$ws.Range("A20").Value = "巨力索具"
$ws.Range("B20").Value = "科大讯飞"
$ws.Range("C20").Value = "工业富联"
$ws.Range("A21").Value = "因赛集团"
$ws.Range("B21").Value = "岩山科技"
$ws.Range("C21").Value = "北方铜业"


Write-Host "Applied HotStock_Top20 cell updates"
